$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tracks")

# --- Update thr_dur (column D) values ---
$ws.Range("D2").Value = 60
$ws.Range("D3").Value = 60
$ws.Range("D4").Value = 48
$ws.Range("D6").Value = 29
$ws.Range("D7").Value = 48

# --- Update calib_2_start (U) / calib_2_end (V) dates for rows 3-6 ---
$ws.Range("U3").Value = 44321
$ws.Range("V3").Value = 44334

$ws.Range("U4").Value = 44321
$ws.Range("V4").Value = 44334

$ws.Range("U5").Value = 44321
$ws.Range("V5").Value = 44334

$ws.Range("U6").Value = 44321
$ws.Range("V6").Value = 44334

# --- Row 7: clear calib_2_start/end/lon/lat (U7:X7) ---
$ws.Range("U7:X7").ClearContents()

# --- Selection / view changes ---
$ws.Range("E8").Select()
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
